# "unify the conception of DataNode, DataTable, Entity."
#
# Rename the two worksheets to their new canonical names:
#   Property1       -> DataNode
#   Record_Station  -> DataTable
$wb = $excel.ActiveWorkbook

$wsDataNode  = $wb.Worksheets.Item("Property1")
$wsDataTable = $wb.Worksheets.Item("Record_Station")

$wsDataNode.Name  = "DataNode"
$wsDataTable.Name = "DataTable"

# "DataTable" becomes the active/selected sheet (previously "DataNode" /
# Property1 held the tab focus), with the selection cursor parked on H32
# instead of the old A10:XFD10 row selection.
$wsDataTable.Activate() | Out-Null
$wsDataTable.Range("H32").Select() | Out-Null
